$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (column E) entirely, shifting
# columns F:K left to E:J.
$ws.Range("E:E").Delete()
